# Applies the OOXML changes described by the diff:
#  1. "Version :" paragraph - merge the two proofErr-wrapped runs into one run.
#  2. "Features :" paragraph - drop the proofErr gramStart/gramEnd markers.
#  3. Insert a new bullet paragraph after the "Refresh only same page..." bullet.
#  4. Move <w:lastRenderedPageBreak/> from the "Ensure proper validation..." run
#     to the "Show a countdown timer..." run.
#  5. Merge the proofErr spellStart/spellEnd-wrapped runs for the block_end_time /
#     attempt_count sentences into single runs (x4 paragraphs).
#  6. Move <w:lastRenderedPageBreak/> from the "Ensure appropriate error
#     messages..." run to the "GET /login-status" run.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Get-ParaByText {
    param($doc, [string]$matchText)
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($matchText)) {
            return $p
        }
    }
    return $null
}

function Set-ParaXml {
    param($doc, [string]$matchText, [string]$innerXml)
    $p = Get-ParaByText $doc $matchText
    if ($null -eq $p) {
        throw "Paragraph containing '$matchText' not found"
    }
    [void]($p.Range.InsertXML($pkgHeader + $innerXml + $pkgFooter))
}

# 1. "Version : 2025.1.1" -- merge the "Version :" + " " runs (proofErr removed)
Set-ParaXml $d "Version : 2025.1.1" ('<w:p w14:paraId="4EEE1646" w14:textId="5C7413A3" w:rsidR="00BC613E" w:rsidRDefault="007477AC">' + `
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Version : </w:t></w:r>' + `
    '<w:r w:rsidR="00000FD4"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2025</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>.</w:t></w:r>' + `
    '<w:r w:rsidR="00000FD4"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>.1</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r>' + `
    '</w:p>')

# 2. "Features :" -- drop proofErr gramStart/gramEnd, keep the two runs as-is
Set-ParaXml $d "Features :" ('<w:p w14:paraId="0DC981FA" w14:textId="02493E77" w:rsidR="00EB0F52" w:rsidRDefault="00C66046">' + `
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00C66046"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Features </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r>' + `
    '</w:p>')

# 3. Insert new bullet paragraph right after the "Refresh only same page..." bullet.
$pRefresh = Get-ParaByText $d "Refresh only same page"
if ($null -eq $pRefresh) {
    throw "Could not find 'Refresh only same page' paragraph"
}
$refreshIndex = $pRefresh.Index
[void]($pRefresh.Range.InsertParagraphAfter())
$newBulletPara = $d.Paragraphs.Item($refreshIndex + 1)
$newBulletXml = $pkgHeader + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>If we enter wrong creds and then write correct creds and enter in the application and again enter wrong creds the login attempts should again update to 3 but its coming 2 login attempts left. (To be decided)</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
[void]($newBulletPara.Range.InsertXML($newBulletXml))

# 4. Move lastRenderedPageBreak from "Ensure proper validation..." to "Show a countdown timer..."
Set-ParaXml $d "Show a countdown timer" ('<w:p w14:paraId="32220CDC" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:lastRenderedPageBreak/><w:t>Show a countdown timer indicating the remaining block time (optional but improves UX).</w:t></w:r>' + `
    '</w:p>')

Set-ParaXml $d "Ensure proper validation" ('<w:p w14:paraId="3B750320" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>Ensure proper validation for incorrect credentials.</w:t></w:r>' + `
    '</w:p>')

# 5. Merge the spellStart/spellEnd proofErr-wrapped runs into single runs.
Set-ParaXml $d "After 3 failed attempts" ('<w:p w14:paraId="7B0F25B6" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="4"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>After 3 failed attempts, set the block_end_time to the current time + 30 minutes.</w:t></w:r>' + `
    '</w:p>')

Set-ParaXml $d "Reset the attempt_count" ('<w:p w14:paraId="5C256D05" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="4"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>Reset the attempt_count and clear block_end_time.</w:t></w:r>' + `
    '</w:p>')

Set-ParaXml $d "is not NULL and greater than the current time" ('<w:p w14:paraId="427FBF01" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>If block_end_time is not NULL and greater than the current time, deny login and provide the block time.</w:t></w:r>' + `
    '</w:p>')

# 6. Move lastRenderedPageBreak from "Ensure appropriate error messages..." to "GET /login-status"
Set-ParaXml $d "GET /login-status" ('<w:p w14:paraId="4A125A87" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>GET /login-status</w:t></w:r>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>: Return the current login attempt count and block status for a user.</w:t></w:r>' + `
    '</w:p>')

Set-ParaXml $d "Ensure appropriate error messages" ('<w:p w14:paraId="3C10C63D" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>Ensure appropriate error messages and HTTP status codes are returned.</w:t></w:r>' + `
    '</w:p>')

# 7. Update database schema sentence -- merge spellStart/spellEnd runs.
Set-ParaXml $d "Update the database schema" ('<w:p w14:paraId="234D47EA" w14:textId="77777777" w:rsidR="00EA2E04" w:rsidRPr="00EA2E04" w:rsidRDefault="00EA2E04" w:rsidP="00EA2E04">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EA2E04"><w:t>Update the database schema to include attempt_count and block_end_time.</w:t></w:r>' + `
    '</w:p>')
